# Updated symbol list on Sun Jan 29 04:55:20 UTC 2023 with GitHub Actions
# Refresh Price (col D) and Volume(1h) (col E) values for the crypto ticker rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.23"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.34%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.64"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.90%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.55%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08128"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.77%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.944"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.87%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.137"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.88%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.25%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9288"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.52%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1422"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.09%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1919"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.77%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09132"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.28%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03515"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.23%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.34%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001391"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.92%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005946"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.13%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.939"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "7.11%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.47%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3429"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.66%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1348"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.43%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.643"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.43%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2450"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.12%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04367"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.40%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.01%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004370"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.70%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.03%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004001"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-10.04%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02042"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.42%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05061"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.22%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007374"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.57%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009835"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.27%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1366"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.24%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002131"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.03%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009373"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.19%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006372"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.86%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.06%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002729"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-18.78%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.06%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.06%"
